$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @{
    2  = @(0.6753301551942219, 1.667794583268128, 3.900430680208489, 8.660232485948974, 0, 14.90378790461981)
    3  = @(0.127881588408715, 0.3127903958511391, 0.8054896365839992, 0.496779210170732, 1, 1.742940831014585)
    4  = @(1.459612070389937, 1.667794583268128, 0.8054896365839992, 0.496779210170732, 0, 4.429675500412797)
    5  = @(3.230985683306322, 1.667794583268128, 0.8054896365839992, 0.496779210170732, 0, 6.201049113329182)
    6  = @(0.6753301551942219, 1.667794583268128, 0.1575252929769615, 0.496779210170732, 0, 2.997429241610044)
    7  = @(1.459612070389937, 0.3127903958511391, 0.8054896365839992, 8.660232485948974, 0, 11.23812458877405)
    8  = @(0.6753301551942219, 1.667794583268128, 0.1575252929769615, 0.496779210170732, 0, 2.997429241610044)
    9  = @(0.3048080303191223, 0.3127903958511391, 0.8054896365839992, 0.496779210170732, 1, 1.919867272924993)
    10 = @(0.3048080303191223, 0.3127903958511391, 0.1575252929769615, 0.496779210170732, 0, 1.271902929317955)
    11 = @(1.459612070389937, 0.04240448674262143, 3.900430680208489, 0.496779210170732, 1, 5.899226447511779)
    12 = @(1.459612070389937, 1.667794583268128, 0.1575252929769615, 0.496779210170732, 1, 3.781711156805759)
    13 = @(1.459612070389937, 1.667794583268128, 0.1575252929769615, 0.496779210170732, 0, 3.781711156805759)
    14 = @(3.230985683306322, 1.667794583268128, 0.8054896365839992, 8.660232485948974, 0, 14.36450238910742)
    15 = @(3.230985683306322, 1.667794583268128, 3.900430680208489, 8.660232485948974, 0, 17.45944343273191)
    16 = @(3.230985683306322, 1.667794583268128, 0.8054896365839992, 0.496779210170732, 1, 6.201049113329182)
    17 = @(3.230985683306322, 1.667794583268128, 0.1575252929769615, 0.496779210170732, 1, 5.553084769722144)
    18 = @(3.230985683306322, 1.667794583268128, 0.1575252929769615, 0.496779210170732, 0, 5.553084769722144)
}

foreach ($row in $data.Keys) {
    $vals = $data[$row]
    $ws.Cells.Item($row, 2).Value = $vals[0]
    $ws.Cells.Item($row, 3).Value = $vals[1]
    $ws.Cells.Item($row, 4).Value = $vals[2]
    $ws.Cells.Item($row, 5).Value = $vals[3]
    $ws.Cells.Item($row, 6).Value = $vals[4]
    $ws.Cells.Item($row, 7).Value = $vals[5]
}
